$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect the new "through" date
$ws.Name = "Through 2022-10-27"

# Update the header label for the current month column (B1)
$ws.Range("B1").Value = "October 2022 (through October 27)"

# Add new carjacking data for 2022-11-04 update
$ws.Range("AF2").Value = 7
$ws.Range("AP2").Value = 8
$ws.Range("BT2").Value = 4
$ws.Range("L3").Value = 9
$ws.Range("AF3").Value = 5
$ws.Range("B5").Value = 2
$ws.Range("L5").Value = 18
$ws.Range("V5").Value = 11
$ws.Range("L6").Value = 13
$ws.Range("V6").Value = 11
$ws.Range("AZ6").Value = 8
$ws.Range("BT6").Value = 4
$ws.Range("B7").Value = 7
$ws.Range("D7").Value = 7
$ws.Range("L7").Value = 6
$ws.Range("AF9").Value = 1
$ws.Range("V11").Value = 7
$ws.Range("AZ11").Value = 3
$ws.Range("V12").Value = 2
$ws.Range("L15").Value = 2
$ws.Range("V16").Value = 5
$ws.Range("AF18").Value = 1
$ws.Range("L23").Value = 5
$ws.Range("AZ23").Value = 1
$ws.Range("AZ24").Value = 3
$ws.Range("V39").Value = 2
$ws.Range("AZ43").Value = 1
$ws.Range("AZ47").Value = 2
$ws.Range("AZ56").Value = 3
$ws.Range("AZ61").Value = 1
$ws.Range("AF98").Value = 1
